$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.193.31'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '2.269.31'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.80'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("E6").Value = '  +3.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.530'
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.30'
$ws.Range("E10").Value = '  +9.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0795'
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.65'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '2.620.16'
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.42'
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").Value = '2.275.76'
$ws.Range("E16").Value = '  +0.47%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.795'
$ws.Range("E17").Value = '  +1.50%  '
$ws.Range("D18").Value = '42.094.56'
$ws.Range("E18").Value = '  +0.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("E19").Value = '  -2.78%  '
$ws.Range("D20").Value = '0.0₃0908'
$ws.Range("E20").Value = '  +0.19%  '
$ws.Range("E21").Value = '  +1.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.02'
$ws.Range("E22").Value = '  +0.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '238.19'
$ws.Range("E23").Value = '  -2.26%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +0.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.73'
$ws.Range("E27").Value = '  -0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.99'
$ws.Range("E28").Value = '  +6.04%  '
$ws.Range("E29").Value = '  -0.42%  '
$ws.Range("E30").Value = '  +2.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.74'
$ws.Range("E31").Value = '  +0.45%  '
$ws.Range("E32").Value = '  -0.76%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.21'
$ws.Range("E34").Value = '  +6.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0741'
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.13'
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.38'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.04'
$ws.Range("E41").Value = '  +2.69%  '
$ws.Range("E42").Value = '  +7.97%  '
$ws.Range("D43").Value = '1.989.08'
$ws.Range("E43").Value = '  -0.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.12'
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.94'
$ws.Range("E46").Value = '  -4.33%  '
$ws.Range("E47").Value = '  +1.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '53.20'
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("E49").Value = '  +0.51%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.18'
$ws.Range("E50").Value = '  -1.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '91.23'
$ws.Range("E51").Value = '  -0.22%  '
